$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three new rows at row 4; this pushes the existing "June 2018" (row 4)
# and "October 2018" (row 5) rows down to rows 7 and 8.
$ws.Rows("4:6").Insert()

# The newly inserted rows inherit blank default formatting; copy the
# border/number-format styling from row 3 (same column layout) down into
# the three new rows so they look like the rest of the table.
$ws.Range("B3:K3").Copy()
$ws.Range("B4:K6").PasteSpecial(-4122)

# New row 4 ("2.4" correlation length / grid group, first sub-row)
$ws.Range("C4").Value = 20
$ws.Range("D4").Value = 2.4
$ws.Range("E4").Value = -158.00619
$ws.Range("F4").Value = 22.738772000000001
$ws.Range("G4").Value = 4729.92
$ws.Range("H4").Value = 0.61
$ws.Range("I4").Value = -0.71699999999999997
$ws.Range("J4").Value = 5.23
$ws.Range("K4").Value = 4735.1499999999996

# New row 5 (second sub-row; only the first three columns populated)
$ws.Range("C5").Value = 8
$ws.Range("D5").Value = 2.4
$ws.Range("E5").Value = -158.00619
$ws.Range("F5").Value = 22.738772000000001
$ws.Range("G5").Value = 4729.92

# New row 6 (third sub-row; only correlation length / grid size populated)
$ws.Range("C6").Value = 20
$ws.Range("D6").Value = 2.4

# Restore the selection the author left behind and refresh the dimension.
$ws.Range("D7").Select()

$wb.Save()
